$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 3: Jacob's peer assessment in "Peer and self assessment" sheet ---
$ws.Range("B3").Value = "Good"

$c3Text = "1)Active at Discord meetings, Good job working on the PoA and sharing`r`nknowledge with other project group members in a scienfiic and respectful`r`nmanner.`r`n2)Quick responsetime on Discord and been doing a good job of managing the weekly discord meetings.  "
$ws.Range("C3").Value = $c3Text
$ws.Range("C3").WrapText = $true

# --- Row 15: Jacob's peer assessment in second block of the same sheet ---
$ws.Range("B15").Value = "Good"

$c15Text = "1) Actively using both Discord and GitHub in an organized way. `r`n2) Great job at sharing knowledge from research and keeping all `r`ngroup members updated by uploading each version of the PoA to GitHub`r`nand explaining changes over Discord "
$ws.Range("C15").Value = $c15Text
$ws.Range("C15").WrapText = $true

# --- View state: scroll position + active selection moved up to rows 13-15 ---
$excel.Goto($ws.Range("B13"), $true) | Out-Null
$ws.Range("C15").Select() | Out-Null
